$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p135v_a1</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p135v_1</id>", 2)
